$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'66.692.26"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +0.36%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'3.487.26"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +0.04%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  -0.01%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'593.35"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +0.72%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'171.91"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +2.60%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  +0.00%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.579"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -1.81%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  +3.24%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  -2.73%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.429"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +0.03%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'4.096.24"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +0.12%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = "'  +0.28%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'29.14"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +4.57%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'66.704.28"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +0.31%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.0000177"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -0.06%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'3.491.50"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -0.01%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'6.24"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -0.07%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'14.27"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +2.63%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'388.62"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -0.33%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'7.90"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +0.54%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'73.16"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +0.69%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'  -0.22%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.532"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +0.64%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'5.67"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -0.83%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'  -0.22%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'10.07"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -0.51%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'  +0.22%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'0.997"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -0.31%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'6.10"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -3.27%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "'  -2.09%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "'  +0.25%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'23.57"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -0.07%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'7.33"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +0.67%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'1.59"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +1.11%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'163.67"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +0.60%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'0.872"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -2.40%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'1.89"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -0.54%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'6.80"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +0.26%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'4.60"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -0.04%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'2.819.86"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +1.91%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'26.91"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +1.98%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.0727"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -0.95%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'25.85"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -1.39%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'42.40"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -0.81%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'2.54"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +0.49%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'  -2.95%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'338.77"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -0.49%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("B49").Value = "'ONDO"
$ws.Range("B49").Style = "Normal"
$ws.Range("C49").Value = "'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("C49").Style = "Normal"
$ws.Range("D49").Value = "'1.07"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -0.40%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("B50").Value = "'Arweave"
$ws.Range("B50").Style = "Normal"
$ws.Range("C50").Value = "'https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("C50").Style = "Normal"
$ws.Range("D50").Value = "'33.90"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +1.82%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("B51").Value = "'Cosmos"
$ws.Range("B51").Style = "Normal"
$ws.Range("C51").Value = "'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("C51").Style = "Normal"
$ws.Range("D51").Value = "'6.40"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -0.91%  "
$ws.Range("E51").Style = "Normal"
